$d = $word.ActiveDocument

# -------------------------------------------------------------------------
# Change 1: append a closing parenthesis after "Gulf Stream Eagle" as its
# own run, with no italic formatting (unlike the preceding "Gulf Stream
# Eagle" run which is italicized).
# -------------------------------------------------------------------------
$rngEagle = $d.Content
$okEagle = $rngEagle.Find.Execute("Gulf Stream Eagle", $true, $false, $false, $false, $false, `
                                   $true, 1, $false, "", 0)
if (-not $okEagle) { throw "Could not find 'Gulf Stream Eagle'" }
$rngEagle.Collapse(0)
$rngEagle.InsertAfter(")")

# -------------------------------------------------------------------------
# Change 2: the "Grossi, M.D.(())" paragraph had stray markdown artifacts
# ("**" prefix and "(())" around "M.D.") - clean the text up to plain
# "Grossi, M.D." and make that name bold (matching the styling used by the
# sibling conference-talk paragraphs), while the remainder of the sentence
# stays a normal, unformatted run.
# -------------------------------------------------------------------------
$old = "**Grossi, M.D.(()), T.M. Özgökmen (2018) Can artificial intelligence predict the dispersion of spilled oil?, Gulf of Mexico Oil Spill and Ecosystem Science Conference, New Orleans, LA."
$new = "Grossi, M.D., T.M. Özgökmen (2018) Can artificial intelligence predict the dispersion of spilled oil?, Gulf of Mexico Oil Spill and Ecosystem Science Conference, New Orleans, LA."

$rngReplace = $d.Content
$okReplace = $rngReplace.Find.Execute($old, $true, $false, $false, $false, $false, `
                                       $true, 1, $false, $new, 2)
if (-not $okReplace) { throw "Could not find/replace the Grossi 2018 paragraph text" }

# Grab the bold "Grossi, M.D." formatting that's already used elsewhere in
# the document (e.g. the very first occurrence), so the fix matches that
# exact styling (w:b + w:bCs).
$rngStyleSource = $d.Content
$okStyle = $rngStyleSource.Find.Execute("Grossi, M.D.", $true, $false, $false, $false, $false, `
                                         $true, 1, $false, "", 0)
if (-not $okStyle) { throw "Could not find a bold 'Grossi, M.D.' style source" }
$boldFormattedText = $rngStyleSource.FormattedText

# Scope the search for our specific "Grossi, M.D." occurrence to the
# stretch of text between the end of the preceding (2019) paragraph and a
# unique marker inside our own (2018) paragraph, so we don't touch any of
# the other "Grossi, M.D." occurrences in the document.
$rngPrev = $d.Content
$okPrev = $rngPrev.Find.Execute("Are Neural Networks Up to the Task?, Gulf of Mexico Oil Spill and Ecosystem Science Conference, New Orleans, LA.", `
                                 $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $okPrev) { throw "Could not find preceding (2019) paragraph" }

$rngAnchor = $d.Content
$okAnchor = $rngAnchor.Find.Execute("Can artificial intelligence predict the dispersion", $true, $false, $false, $false, $false, `
                                     $true, 1, $false, "", 0)
if (-not $okAnchor) { throw "Could not find anchor text inside the (2018) paragraph" }

$rngTarget = $d.Range($rngPrev.End, $rngAnchor.Start)
$okTarget = $rngTarget.Find.Execute("Grossi, M.D.", $true, $false, $false, $false, $false, `
                                     $true, 1, $false, "", 0)
if (-not $okTarget) { throw "Could not find the 'Grossi, M.D.' run to bold" }

$rngTarget.FormattedText = $boldFormattedText

Write-Host "Done."
